$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 2) ---
# J2 is set up first so its new "wrap text / top aligned" style becomes the
# next style table entry, matching the added cellXfs entry for the HSN code
# list cell.
$ws.Range("J2").Value = "996211`n62052000`n62052000`n62046200`n48194000`n33072000`n39264099"
$ws.Range("J2").VerticalAlignment = -4160
$ws.Range("J2").WrapText = $true

# Plain text cells - Excel keeps these as literal text already.
$ws.Range("A2").Value = "WESTSIDE`nSjr Zion, Survey"
$ws.Range("B2").Value = "29AAACL1838J1ZC"
$ws.Range("C2").Value = "W089100169940"
$ws.Range("I2").Value = "N/A"

# Columns D-H hold values that look like a date / numbers ("2024-09-28",
# "4045.01", "173.91" ...) but must be stored as literal text, matching the
# source export (inlineStr cells in the target workbook). Temporarily force
# text formatting so Excel's COM layer doesn't auto-convert them into a
# date serial / numeric value, write the literal strings, then restore the
# cells' original (General / unformatted) appearance by pasting the format
# back in from an untouched, plain cell - this avoids leaving any visible
# formatting on these cells, same as the rest of the row.
$target = $ws.Range("D2:H2")
$target.NumberFormat = "@"
$ws.Range("D2").Value = "2024-09-28"
$ws.Range("E2").Value = "4045.01"
$ws.Range("F2").Value = "173.91"
$ws.Range("G2").Value = "173.91"
$ws.Range("H2").Value = "173.91"

$ws.Range("D5:H5").Copy()
$target.PasteSpecial(-4122)
$excel.CutCopyMode = $false
